# Generate Report for Handoff
#
# Refreshes the localization-status report for the rows whose handoff
# xliff was (re)generated: status "Ready for handoff" rows 7, 8, 10, 11,
# 12, 14 on the zh-cn and de-de sheets now have their Priority set to
# "ht" (handoff type), and the various "latest handoff generated" /
# "latest handoff datetime" timestamps for those rows move forward ~32s
# to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = 7,8,10,11,12,14

foreach ($r in $rows) {
    # Priority column (E) on both locale sheets: blank -> "ht"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Latest Handoff Datetime (H) per locale sheet
    $zhcn.Range("H$r").Value = "2016-08-12 22:26:08"
    $dede.Range("H$r").Value = "2016-08-12 22:26:16"

    # Latest HO Xliff Generate Date (G) on the Overview rollup sheet
    $overview.Range("G$r").Value = "2016-08-12 22:26:16"
}
